$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the "phone" column (B). Excel shifts service_user (C) into B and
#    firm_ref (D) into C, and removes the now-orphaned phone strings.
# ---------------------------------------------------------------------------
$ws.Columns("B").Delete()

# ---------------------------------------------------------------------------
# 2. Seed data changes.
#    - admin@wildwheeladventures.cz -> admin@stopoverx.com, now a service user (TRUE)
#    - admin@glidequesttours.cz row is cleared out (firm admin removed from seed)
#    - firm_ref column (C) data is no longer seeded, only the header stays
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "admin@stopoverx.com"
$ws.Range("B2").Value = $true

$ws.Range("A3").ClearContents()

$ws.Range("C2:C10").Clear()

# ---------------------------------------------------------------------------
# 3. Rebuild the hyperlinks: drop the old ones (wildwheeladventures.cz /
#    glidequesttours.cz admins) and re-add the rest in row order so the
#    relationship ids line up the way Excel would renumber them.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:admin@stopoverx.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:serviceuser-wildwheeladventures@stopover.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:serviceuser-glidequesttours@stopover.com")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:customer1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:customer2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:customer3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:customer4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:mikhail@dorokhovich.ru")

# A3 no longer links anywhere, keep its plain (non-hyperlink) look but restore
# the same "email style" formatting the rest of column A uses.
$ws.Range("A3").Font.Underline = $true
$ws.Range("A3").Font.Name = "Calibri (Body)"

# A2 becomes a real/plain hyperlink -> use the true Excel "Hyperlink" style.
$ws.Range("A2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 4. Selection moves to A2 after the edit.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()

"done"
